# Auto-generated Excel COM-interop script
# Applies scheduled market-price refresh updates to the Balmung_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1089.52
$ws.Range("I18").Value = 1089.52
$ws.Range("K18").Value = 1089.52
$ws.Range("M18").Value = -805.52
$ws.Range("H33").Value = 3793415.8
$ws.Range("I33").Value = 5130919
$ws.Range("K33").Value = 5130919
$ws.Range("M33").Value = -5130690
$ws.Range("H53").Value = 250000940
$ws.Range("I53").Value = 300
$ws.Range("K53").Value = 300
$ws.Range("M53").Value = 337
$ws.Range("H62").Value = 2590.5
$ws.Range("I62").Value = 2253.4285
$ws.Range("K62").Value = 2253.4285
$ws.Range("M62").Value = -1629.4285
$ws.Range("H65").Value = 2590.5
$ws.Range("I65").Value = 2253.4285
$ws.Range("K65").Value = 11267.1425
$ws.Range("M65").Value = -8147.1425
$ws.Range("H69").Value = 14296623
$ws.Range("J69").Value = 12880
$ws.Range("L69").Value = 38640
$ws.Range("N69").Value = -40388
$ws.Range("H72").Value = 14296623
$ws.Range("J72").Value = 12880
$ws.Range("L72").Value = 115920
$ws.Range("N72").Value = -124656
$ws.Range("H80").Value = 10869988
$ws.Range("I80").Value = 301.27274
$ws.Range("J80").Value = 20833868
$ws.Range("K80").Value = 903.81822
$ws.Range("L80").Value = 62501604
$ws.Range("M80").Value = 94.18178
$ws.Range("N80").Value = -62503600
$ws.Range("H83").Value = 10869988
$ws.Range("I83").Value = 301.27274
$ws.Range("J83").Value = 20833868
$ws.Range("K83").Value = 2711.45466
$ws.Range("L83").Value = 187504812
$ws.Range("M83").Value = 2280.54534
$ws.Range("N83").Value = -187514796
$ws.Range("H88").Value = 10636.333
$ws.Range("I88").Value = 1399
$ws.Range("K88").Value = 1399
$ws.Range("M88").Value = -993
$ws.Range("H91").Value = 10636.333
$ws.Range("I91").Value = 1399
$ws.Range("K91").Value = 1399
$ws.Range("M91").Value = 5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 879.8570999999999
$ws.Range("H45").Value = 43077.24
$ws.Range("I45").Value = 64139.625
$ws.Range("K45").Value = 64139.625
$ws.Range("M45").Value = -63762.625
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H88").Value = 2031.9474
$ws.Range("I88").Value = 1900
$ws.Range("J88").Value = 2047.4706
$ws.Range("K88").Value = 1900
$ws.Range("L88").Value = 2047.4706
$ws.Range("M88").Value = -1494
$ws.Range("N88").Value = -2859.4706
$ws.Range("H91").Value = 2031.9474
$ws.Range("I91").Value = 1900
$ws.Range("J91").Value = 2047.4706
$ws.Range("K91").Value = 1900
$ws.Range("L91").Value = 2047.4706
$ws.Range("M91").Value = -496
$ws.Range("N91").Value = -4855.470600000001
$ws.Range("H97").Value = 4925.7407
$ws.Range("I97").Value = 5511.136
$ws.Range("K97").Value = 5511.136
$ws.Range("M97").Value = -5015.136
$ws.Range("H112").Value = 49795.332
$ws.Range("J112").Value = 49795.332
$ws.Range("L112").Value = 49795.332
$ws.Range("N112").Value = -52749.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 879.8570999999999
$ws.Range("H86").Value = 8160.3125
$ws.Range("I86").Value = 5869
$ws.Range("J86").Value = 13201.2
$ws.Range("K86").Value = 5869
$ws.Range("L86").Value = 13201.2
$ws.Range("M86").Value = -4746
$ws.Range("N86").Value = -15447.2
$ws.Range("H89").Value = 8160.3125
$ws.Range("I89").Value = 5869
$ws.Range("J89").Value = 13201.2
$ws.Range("K89").Value = 29345
$ws.Range("L89").Value = 66006
$ws.Range("M89").Value = -23729
$ws.Range("N89").Value = -77238
$ws.Range("H105").Value = 15444.308
$ws.Range("I105").Value = 12148
$ws.Range("K105").Value = 12148
$ws.Range("M105").Value = -10401
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 62
$ws.Range("I7").Value = 61.666668
$ws.Range("J7").Value = 62.5
$ws.Range("K7").Value = 61.666668
$ws.Range("L7").Value = 62.5
$ws.Range("M7").Value = 51.333332
$ws.Range("N7").Value = -288.5
$ws.Range("H10").Value = 3992
$ws.Range("J10").Value = 4234.5
$ws.Range("L10").Value = 4234.5
$ws.Range("N10").Value = -4512.5
$ws.Range("H22").Value = 1860.3334
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H132").Value = 1844.5652
$ws.Range("I132").Value = 1809.375
$ws.Range("J132").Value = 1925
$ws.Range("K132").Value = 5428.125
$ws.Range("L132").Value = 5775
$ws.Range("M132").Value = -2898.125
$ws.Range("N132").Value = -10835

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1271.125
$ws.Range("J2").Value = 1719.95
$ws.Range("L2").Value = 10319.7
$ws.Range("N2").Value = -10545.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 31261546
$ws.Range("I80").Value = 9035
$ws.Range("J80").Value = 55569056
$ws.Range("K80").Value = 9035
$ws.Range("L80").Value = 55569056
$ws.Range("M80").Value = -8037
$ws.Range("N80").Value = -55571052
$ws.Range("H83").Value = 31261546
$ws.Range("I83").Value = 9035
$ws.Range("J83").Value = 55569056
$ws.Range("K83").Value = 45175
$ws.Range("L83").Value = 277845280
$ws.Range("M83").Value = -40183
$ws.Range("N83").Value = -277855264

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8213.333000000001
$ws.Range("I46").Value = 35999.332
$ws.Range("J46").Value = 3582.3333
$ws.Range("K46").Value = 35999.332
$ws.Range("L46").Value = 3582.3333
$ws.Range("M46").Value = -35811.332
$ws.Range("N46").Value = -3958.3333
$ws.Range("H68").Value = 15000
$ws.Range("I68").Value = 15000
$ws.Range("K68").Value = 15000
$ws.Range("M68").Value = -14251
$ws.Range("H71").Value = 15000
$ws.Range("I71").Value = 15000
$ws.Range("K71").Value = 75000
$ws.Range("M71").Value = -71256
$ws.Range("H93").Value = 1237.2222
$ws.Range("I93").Value = 1237.2222
$ws.Range("K93").Value = 1237.2222
$ws.Range("M93").Value = 10.77780000000007
$ws.Range("H139").Value = 189990
$ws.Range("J139").Value = 189990
$ws.Range("L139").Value = 189990
$ws.Range("N139").Value = -200270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 60028.824
$ws.Range("I81").Value = 999.625
$ws.Range("J81").Value = 112499.22
$ws.Range("K81").Value = 1999.25
$ws.Range("L81").Value = 224998.44
$ws.Range("M81").Value = -938.25
$ws.Range("N81").Value = -227120.44
$ws.Range("H84").Value = 60028.824
$ws.Range("I84").Value = 999.625
$ws.Range("J84").Value = 112499.22
$ws.Range("K84").Value = 9996.25
$ws.Range("L84").Value = 1124992.2
$ws.Range("M84").Value = -4692.25
$ws.Range("N84").Value = -1135600.2
$ws.Range("H132").Value = 1923.7838
$ws.Range("I132").Value = 1491.6428
$ws.Range("K132").Value = 4474.928400000001
$ws.Range("M132").Value = -1944.928400000001
